# Adder - BOM update
# Move quantities from the "Bought" column (D) to the "Have" column (C)
# for rows 3, 12, 15 and 16, matching the rest of the BOM rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was D3=20 -> now C3=2 (D3 cleared)
$ws.Range("C3").Value = 2
$ws.Range("D3").ClearContents()

# Row 12: was C12=4, D12=10 -> now C12=12 (D12 cleared)
$ws.Range("C12").Value = 12
$ws.Range("D12").ClearContents()

# Row 15: was D15=2 -> now C15=2 (D15 cleared)
$ws.Range("C15").Value = 2
$ws.Range("D15").ClearContents()

# Row 16: was D16=2 -> now C16=1 (D16 cleared)
$ws.Range("C16").Value = 1
$ws.Range("D16").ClearContents()

# Update the active selection to reflect where the user left off editing
$ws.Range("C17").Select() | Out-Null
